# demeaned rs scores (binning)
# Recompute Composite_Reward (col B) as the mean-centered ("demeaned") raw
# reward score, and Composite_Reward_Squared (col C) as the square of the
# demeaned score.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mean = 5.52

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 51 }

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $rawValue = $bCell.Value2
    if ($rawValue -eq $null) { continue }

    $demeaned = $rawValue - $mean
    $squared = $demeaned * $demeaned

    $bCell.Value = $demeaned
    $ws.Cells.Item($r, 3).Value = $squared
}
